$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rooms")

# Row 2 values (forced text via NumberFormat, reset style afterward to keep default style index)
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "100"
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2023-01-10"
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "2023-01-20"
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "2023-01-10"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2023-01-20"
$ws.Range("E2").Style = "Normal"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "2023-01-10"
$ws.Range("F2").Style = "Normal"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "2023-01-20"
$ws.Range("G2").Style = "Normal"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "2024-01-10"
$ws.Range("H2").Style = "Normal"
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "2024-01-15"
$ws.Range("I2").Style = "Normal"
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "2024-01-10"
$ws.Range("J2").Style = "Normal"
$ws.Range("K2").NumberFormat = "@"
$ws.Range("K2").Value = "2024-01-15"
$ws.Range("K2").Style = "Normal"
$ws.Range("L2").NumberFormat = "@"
$ws.Range("L2").Value = "2024-01-10"
$ws.Range("L2").Style = "Normal"
$ws.Range("M2").NumberFormat = "@"
$ws.Range("M2").Value = "2024-01-15"
$ws.Range("M2").Style = "Normal"
$ws.Range("N2").NumberFormat = "@"
$ws.Range("N2").Value = "2024-01-20"
$ws.Range("N2").Style = "Normal"
$ws.Range("O2").NumberFormat = "@"
$ws.Range("O2").Value = "2024-01-30"
$ws.Range("O2").Style = "Normal"
$ws.Range("P2").NumberFormat = "@"
$ws.Range("P2").Value = "2024-01-20"
$ws.Range("P2").Style = "Normal"
$ws.Range("Q2").NumberFormat = "@"
$ws.Range("Q2").Value = "2024-01-30"
$ws.Range("Q2").Style = "Normal"
$ws.Range("R2").NumberFormat = "@"
$ws.Range("R2").Value = "2024-01-20"
$ws.Range("R2").Style = "Normal"
$ws.Range("S2").NumberFormat = "@"
$ws.Range("S2").Value = "2024-01-30"
$ws.Range("S2").Style = "Normal"
$ws.Range("T2").NumberFormat = "@"
$ws.Range("T2").Value = "2024-02-20"
$ws.Range("T2").Style = "Normal"
$ws.Range("U2").NumberFormat = "@"
$ws.Range("U2").Value = "2024-02-28"
$ws.Range("U2").Style = "Normal"
$ws.Range("V2").NumberFormat = "@"
$ws.Range("V2").Value = "2024-02-20"
$ws.Range("V2").Style = "Normal"
$ws.Range("W2").NumberFormat = "@"
$ws.Range("W2").Value = "2024-02-28"
$ws.Range("W2").Style = "Normal"
$ws.Range("X2").NumberFormat = "@"
$ws.Range("X2").Value = "2024-02-20"
$ws.Range("X2").Style = "Normal"
$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value = "2024-02-28"
$ws.Range("Y2").Style = "Normal"

# Empty trailing cells Z2:AU2 (present in sheetData but blank)
$ws.Range("Z2:AU2").NumberFormat = "@"
$ws.Range("Z2:AU2").Style = "Normal"

# Room numbers column A, rows 3-5
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "200"
$ws.Range("A3").Style = "Normal"
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "300"
$ws.Range("A4").Style = "Normal"
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "400"
$ws.Range("A5").Style = "Normal"
